# Fruta / hortaliza, semanal
# A new weekly price report (week of 2023-07-24, serial 45131) was added for
# "Agrícola del Norte S.A. de Arica - Frutilla". The new week's 4 quality
# rows (Especial/Primera/Segunda/Tercera) are inserted at the top of the
# existing data block (rows 56-59), pushing the previous weeks down by 4
# rows (old rows 56-68 become rows 60-72).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before row 56, shifting the rest of the table down.
$ws.Rows("56:59").Insert()

$newWeek = @(
  @{ Row = 56; Calidad = "Especial"; Volumen = 370; PrecioMin = 7000; PrecioMax = 8000; PrecioProm = 7405; PrecioKg = 2468 },
  @{ Row = 57; Calidad = "Primera";  Volumen = 370; PrecioMin = 5000; PrecioMax = 6000; PrecioProm = 5514; PrecioKg = 1838 },
  @{ Row = 58; Calidad = "Segunda";  Volumen = 490; PrecioMin = 4000; PrecioMax = 5000; PrecioProm = 4449; PrecioKg = 1483 },
  @{ Row = 59; Calidad = "Tercera";  Volumen = 470; PrecioMin = 3000; PrecioMax = 4000; PrecioProm = 3426; PrecioKg = 1142 }
)

foreach ($item in $newWeek) {
  $r = $item.Row
  $ws.Cells.Item($r, 1).Value = 1
  $ws.Cells.Item($r, 2).Value = "Agrícola del Norte S.A. de Arica"
  $ws.Cells.Item($r, 3).Value = "Arica y Parinacota"
  $ws.Cells.Item($r, 4).Value = 45131
  $ws.Cells.Item($r, 5).Value = 15
  $ws.Cells.Item($r, 6).Value = "Fruta"
  $ws.Cells.Item($r, 7).Value = 100101
  $ws.Cells.Item($r, 8).Value = "Berries"
  $ws.Cells.Item($r, 9).Value = 100112025
  $ws.Cells.Item($r, 10).Value = "Frutilla"
  $ws.Cells.Item($r, 11).Value = "Sin especificar"
  $ws.Cells.Item($r, 12).Value = $item.Calidad
  $ws.Cells.Item($r, 13).Value = $item.Volumen
  $ws.Cells.Item($r, 14).Value = $item.PrecioMin
  $ws.Cells.Item($r, 15).Value = $item.PrecioMax
  $ws.Cells.Item($r, 16).Value = $item.PrecioProm
  $ws.Cells.Item($r, 17).Value = "`$/bandeja 3 kilos"
  $ws.Cells.Item($r, 18).Value = "Región de Arica y Parinacota"
  $ws.Cells.Item($r, 19).Value = $item.PrecioKg
  $ws.Cells.Item($r, 20).Value = 3
}
